# Disaggregation of commodity Copper
# Rename the commodity label "Copper ores and concentrates" -> "Copper"
# on every year sheet (the label lives in row 4, column C of each sheet),
# and apply the small recalculated value corrections to column D row 4
# on the handful of sheets where the totals shifted.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("C4")
    if ($cell.Value2 -eq "Copper ores and concentrates") {
        $cell.Value = "Copper"
    }
}

$corrections = @{
    "2026" = 31910.8392583942
    "2030" = 53707.83886834714
    "2041" = 192438.6403801433
    "2048" = 659989.6169505299
    "2065" = 744340.2861112709
    "2073" = 719459.3854483
    "2074" = 848045.1390536642
}

foreach ($year in $corrections.Keys) {
    $yearName = [string]$year
    $ws = $wb.Worksheets.Item($yearName)
    $ws.Range("D4").Value = $corrections[$year]
}
